# Auto update Excel log
# Appends newly-captured sensor readings (2026-02-06, ~10:08-10:09) to the
# PIR, Humidity and Temperature logs. All cells in these sheets are stored
# as literal text (dates/times/readings as strings), so every value is
# forced to Text format before being written and the explicit format is
# cleared again afterwards so the new cells match the existing plain,
# un-styled cells in the sheet.

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

function Add-LogRows($ws, $rows) {
    foreach ($row in $rows) {
        $r = $row[0]
        for ($c = 1; $c -le 6; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            Set-TextCell $cell $row[$c]
        }
    }
}

$wb = $excel.ActiveWorkbook

# --- PIR sheet: new rows 340-352 -------------------------------------------------
$pirData = @(
    ,@(340, '2026-02-06', '10:08:11', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(341, '2026-02-06', '10:08:12', '10:00', 'Bathroom', 'Motion Detected', 'Active')
    ,@(342, '2026-02-06', '10:08:17', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(343, '2026-02-06', '10:08:22', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(344, '2026-02-06', '10:08:27', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(345, '2026-02-06', '10:08:33', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(346, '2026-02-06', '10:08:37', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(347, '2026-02-06', '10:08:43', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(348, '2026-02-06', '10:08:48', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(349, '2026-02-06', '10:08:53', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(350, '2026-02-06', '10:08:58', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(351, '2026-02-06', '10:09:03', '10:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(352, '2026-02-06', '10:09:08', '10:00', 'Bathroom', 'No Motion', 'Inactive')
)

$wsPIR = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPIR $pirData

# --- Humidity sheet: new rows 228-235 ---------------------------------------------
$humidityData = @(
    ,@(228, '2026-02-06', '10:08:13', '10:00', 'Bathroom', '69.2%', 'Active')
    ,@(229, '2026-02-06', '10:08:32', '10:00', 'Bathroom', '69.2%', 'Active')
    ,@(230, '2026-02-06', '10:08:42', '10:00', 'Bathroom', '69.3%', 'Active')
    ,@(231, '2026-02-06', '10:08:47', '10:00', 'Bathroom', '69.2%', 'Active')
    ,@(232, '2026-02-06', '10:08:52', '10:00', 'Bathroom', '69.2%', 'Active')
    ,@(233, '2026-02-06', '10:08:57', '10:00', 'Bathroom', '69.2%', 'Active')
    ,@(234, '2026-02-06', '10:09:02', '10:00', 'Bathroom', '69.0%', 'Active')
    ,@(235, '2026-02-06', '10:09:07', '10:00', 'Bathroom', '68.9%', 'Active')
)

$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity $humidityData

# --- Temperature sheet: new rows 228-235 ------------------------------------------
$temperatureData = @(
    ,@(228, '2026-02-06', '10:08:13', '10:00', 'Bathroom', '28.0C', 'Active')
    ,@(229, '2026-02-06', '10:08:32', '10:00', 'Bathroom', '28.0C', 'Active')
    ,@(230, '2026-02-06', '10:08:42', '10:00', 'Bathroom', '28.1C', 'Active')
    ,@(231, '2026-02-06', '10:08:48', '10:00', 'Bathroom', '28.1C', 'Active')
    ,@(232, '2026-02-06', '10:08:52', '10:00', 'Bathroom', '28.2C', 'Active')
    ,@(233, '2026-02-06', '10:08:58', '10:00', 'Bathroom', '28.1C', 'Active')
    ,@(234, '2026-02-06', '10:09:02', '10:00', 'Bathroom', '28.1C', 'Active')
    ,@(235, '2026-02-06', '10:09:08', '10:00', 'Bathroom', '28.1C', 'Active')
)

$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature $temperatureData
